$wb = $excel.ActiveWorkbook

$notes = $wb.Worksheets.Item("Notes")
$data = $wb.Worksheets.Item("Data")

# Fix the "Units of measure" note on the Notes sheet.
$notes.Range("A3").Value = "Units of measure: constant 2015 US$"

# Populate the user data row on the Data sheet.
$data.Range("A2").Value = "bilateral-unspecified"
$data.Range("B2").Value = "Bilateral, unspecified"
$data.Range("C2").Value = 2015
$data.Range("D2").Value = 19051100
